$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 632.44446
$ws.Range("I28").Value = 520.875
$ws.Range("J28").Value = 1525
$ws.Range("K28").Value = 520.875
$ws.Range("L28").Value = 1525
$ws.Range("M28").Value = -35.875
$ws.Range("N28").Value = -2495
$ws.Range("H132").Value = 6707856.5
$ws.Range("I132").Value = 7961419.5
$ws.Range("K132").Value = 23884258.5
$ws.Range("M132").Value = -23881728.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 19726
$ws.Range("I45").Value = 14429.25
$ws.Range("K45").Value = 14429.25
$ws.Range("M45").Value = -14052.25
$ws.Range("H102").Value = 3575.25
$ws.Range("I102").Value = 3575.25
$ws.Range("K102").Value = 3575.25
$ws.Range("M102").Value = -1953.25
$ws.Range("H122").Value = 3183.4583
$ws.Range("I122").Value = 3038.4211
$ws.Range("J122").Value = 3734.6
$ws.Range("K122").Value = 9115.263300000001
$ws.Range("L122").Value = 11203.8
$ws.Range("M122").Value = -6665.263300000001
$ws.Range("N122").Value = -16103.8
$ws.Range("H132").Value = 2544.0962
$ws.Range("I132").Value = 2182.0952
$ws.Range("K132").Value = 6546.285600000001
$ws.Range("M132").Value = -4016.285600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3678.9333
$ws.Range("I86").Value = 2579.9473
$ws.Range("J86").Value = 5577.1816
$ws.Range("K86").Value = 2579.9473
$ws.Range("L86").Value = 5577.1816
$ws.Range("M86").Value = -1456.9473
$ws.Range("N86").Value = -7823.1816
$ws.Range("H89").Value = 3678.9333
$ws.Range("I89").Value = 2579.9473
$ws.Range("J89").Value = 5577.1816
$ws.Range("K89").Value = 12899.7365
$ws.Range("L89").Value = 27885.908
$ws.Range("M89").Value = -7283.736499999999
$ws.Range("N89").Value = -39117.908
$ws.Range("H94").Value = 961.0454999999999
$ws.Range("I94").Value = 974.15
$ws.Range("J94").Value = 830
$ws.Range("K94").Value = 974.15
$ws.Range("L94").Value = 830
$ws.Range("M94").Value = -523.15
$ws.Range("N94").Value = -1732
$ws.Range("H107").Value = 3230.2222
$ws.Range("I107").Value = 3311.7144
$ws.Range("K107").Value = 3311.7144
$ws.Range("M107").Value = -1391.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2604.1333
$ws.Range("I16").Value = 2218.889
$ws.Range("J16").Value = 3182
$ws.Range("K16").Value = 2218.889
$ws.Range("L16").Value = 3182
$ws.Range("M16").Value = -1931.889
$ws.Range("N16").Value = -3756
$ws.Range("H105").Value = 1433
$ws.Range("I105").Value = 1433
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1433
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 314
$ws.Range("H109").Value = 35445
$ws.Range("J109").Value = 35445
$ws.Range("L109").Value = 35445
$ws.Range("N109").Value = -37525
$ws.Range("H113").Value = 2604.1333
$ws.Range("I113").Value = 2218.889
$ws.Range("J113").Value = 3182
$ws.Range("K113").Value = 2218.889
$ws.Range("L113").Value = 3182
$ws.Range("M113").Value = -48.88900000000012
$ws.Range("N113").Value = -7522
$ws.Range("H122").Value = 1544.5454
$ws.Range("I122").Value = 1619
$ws.Range("K122").Value = 4857
$ws.Range("M122").Value = -2407
$ws.Range("H141").Value = 51820
$ws.Range("J141").Value = 51820
$ws.Range("L141").Value = 51820
$ws.Range("N141").Value = -62180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2779133.5
$ws.Range("J68").Value = 3572887.2
$ws.Range("L68").Value = 10718661.6
$ws.Range("N68").Value = -10720283.6
$ws.Range("H71").Value = 2779133.5
$ws.Range("J71").Value = 3572887.2
$ws.Range("L71").Value = 32155984.8
$ws.Range("N71").Value = -32164096.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 250000
$ws.Range("J42").Value = 250000
$ws.Range("L42").Value = 250000
$ws.Range("N42").Value = -250970
$ws.Range("H51").Value = 99990
$ws.Range("J51").Value = 99990
$ws.Range("L51").Value = 99990
$ws.Range("N51").Value = -101008
$ws.Range("H102").Value = 2441.4167
$ws.Range("I102").Value = 2163.3635
$ws.Range("K102").Value = 2163.3635
$ws.Range("M102").Value = -541.3634999999999
$ws.Range("H113").Value = 3759.1904
$ws.Range("I113").Value = 3082.1428
$ws.Range("K113").Value = 3082.1428
$ws.Range("M113").Value = -912.1428000000001
$ws.Range("H115").Value = 250000
$ws.Range("J115").Value = 250000
$ws.Range("L115").Value = 250000
$ws.Range("N115").Value = -252350
$ws.Range("H122").Value = 21420.615
$ws.Range("I122").Value = 22465
$ws.Range("K122").Value = 67395
$ws.Range("M122").Value = -64945
$ws.Range("H126").Value = 18747.4
$ws.Range("I126").Value = 22996.666
$ws.Range("K126").Value = 68989.99800000001
$ws.Range("M126").Value = -66519.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13359.3
$ws.Range("I7").Value = 18033.334
$ws.Range("J7").Value = 6348.25
$ws.Range("K7").Value = 18033.334
$ws.Range("L7").Value = 6348.25
$ws.Range("M7").Value = -17921.334
$ws.Range("N7").Value = -6572.25
$ws.Range("H55").Value = 1743.4445
$ws.Range("I55").Value = 1743.4445
$ws.Range("K55").Value = 1743.4445
$ws.Range("M55").Value = -1570.4445
$ws.Range("H122").Value = 7284.5713
$ws.Range("I122").Value = 6198.4
$ws.Range("K122").Value = 18595.2
$ws.Range("M122").Value = -16145.2
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H126").Value = 13359.3
$ws.Range("I126").Value = 18033.334
$ws.Range("J126").Value = 6348.25
$ws.Range("K126").Value = 54100.00199999999
$ws.Range("L126").Value = 19044.75
$ws.Range("M126").Value = -51630.00199999999
$ws.Range("N126").Value = -23984.75
$ws.Range("H132").Value = 5728.5713
$ws.Range("I132").Value = 6000
$ws.Range("K132").Value = 18000
$ws.Range("M132").Value = -15470
$ws.Range("H136").Value = 4067.5625
$ws.Range("I136").Value = 3248.8
$ws.Range("K136").Value = 9746.400000000001
$ws.Range("M136").Value = -7196.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6054.56
$ws.Range("I122").Value = 2233.261
$ws.Range("K122").Value = 6699.782999999999
$ws.Range("M122").Value = -4249.782999999999
$ws.Range("H126").Value = 4502
$ws.Range("I126").Value = 3702.4
$ws.Range("J126").Value = 8500
$ws.Range("K126").Value = 11107.2
$ws.Range("L126").Value = 25500
$ws.Range("M126").Value = -8637.200000000001
$ws.Range("N126").Value = -30440
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920
$ws.Range("H132").Value = 5346.2383
$ws.Range("I132").Value = 6081
$ws.Range("K132").Value = 18243
$ws.Range("M132").Value = -15713
$ws.Range("H136").Value = 1813.9584
$ws.Range("I136").Value = 1319.5
$ws.Range("J136").Value = 3297.3333
$ws.Range("K136").Value = 3958.5
$ws.Range("L136").Value = 9891.999899999999
$ws.Range("M136").Value = -1408.5
$ws.Range("N136").Value = -14991.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N105").ClearContents()
